$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style used by the other
# header cells (copy style from H1 so formatting - bold, border, centered -
# matches the existing headers).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the I0 / IF data values for rows 2-8
$data = @(
    @(1, 2),
    @(8, 10),
    @(6, 6),
    @(4, 5),
    @(9, 9),
    @(1, 2),
    @(1, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]  # column I
    $ws.Cells.Item($row, 10).Value = $data[$i][1] # column J
}
